# SCD0016-043 update: rename sheet, update TC_ID, relabel test case id,
# left-align the data table and move the selection to B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet (was SCD0269)
$ws.Name = "SCD0016"

# 2. Update the TC_ID column (B2:B6) from "DGS-284" to "SCD0016-043"
$ws.Range("B2:B6").Value = "SCD0016-043"

# 3. Apply left + vertically-centered alignment across the header + data rows (A1:P6)
$rng = $ws.Range("A1:P6")
$rng.VerticalAlignment = -4108
$rng.HorizontalAlignment = -4131

# 4. Move the active selection to B7 (matches the saved cursor position)
$ws.Range("B7").Select() | Out-Null
